$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.561.78'
$ws.Range('D3').Value = '1.741.18'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.24'
$ws.Range('E5').Value = '  +5.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4807'
$ws.Range('E7').Value = '  +3.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2690'
$ws.Range('E8').Value = '  +4.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06255'
$ws.Range('E9').Value = '  +1.82%  '
$ws.Range('D10').Value = '1.743.29'
$ws.Range('E10').Value = '  +4.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07139'
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.79'
$ws.Range('E12').Value = '  +7.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6213'
$ws.Range('E13').Value = '  +8.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.521'
$ws.Range('E14').Value = '  +4.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.38'
$ws.Range('E15').Value = '  +3.27%  '
$ws.Range('D17').Value = '26.567.42'
$ws.Range('E17').Value = '  +4.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006902'
$ws.Range('E19').Value = '  +3.12%  '
$ws.Range('E20').Value = '  +3.58%  '
$ws.Range('D21').Value = '1.967.06'
$ws.Range('E21').Value = '  +4.53%  '
$ws.Range('E22').Value = '  +4.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.876'
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.368'
$ws.Range('E24').Value = '  +3.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.04'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.39'
$ws.Range('E26').Value = '  +3.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.813'
$ws.Range('E27').Value = '  +6.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.421'
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '107.05'
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.741'
$ws.Range('E31').Value = '  +3.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07885'
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04576'
$ws.Range('E33').Value = '  +6.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.626'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6380'
$ws.Range('E35').Value = '  +6.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +5.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9379'
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '113.35'
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.432'
$ws.Range('E39').Value = '  -3.19%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.980'
$ws.Range('E40').Value = '  +8.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.756'
$ws.Range('E42').Value = '  +17.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.01514'
$ws.Range('E43').Value = '  +3.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3914'
$ws.Range('E44').Value = '  +5.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.793'
$ws.Range('E45').Value = '  +10.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1212'
$ws.Range('E46').Value = '  +9.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05332'
$ws.Range('E47').Value = '  +1.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.928'
$ws.Range('E48').Value = '  +7.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.76'
$ws.Range('E49').Value = '  +3.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.262'
$ws.Range('E50').Value = '  +5.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3447'
$ws.Range('E51').Value = '  +4.60%  '
